# Update crypto price/volume/hora data as per the Tue Jan 17 04:11:25 UTC 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "299.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.60%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "4"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.99%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "4"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.099"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.39%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "4"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07936"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "6.17%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "4"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.234"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-7.22%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "4"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.755"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.25%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "4"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.858"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.28%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "4"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9155"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.13%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "4"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1732"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.12%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "4"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07370"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.96%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "4"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09423"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "14.70%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "4"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03018"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.23%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "4"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1003"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.02%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "4"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001507"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.03%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "4"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005848"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.63%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "4"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.477"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.84%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "4"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.32%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "4"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.30%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "4"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.19%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "4"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.911"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-15.94%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "4"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1699"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.50%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "4"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04625"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.52%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "4"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.56%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "4"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004466"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.49%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "4"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001199"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.68%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "4"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003395"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "23.92%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "4"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "4"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "4"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "4"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "4"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "4"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "4"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "4"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "4"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "4"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "4"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "4"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01750"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.76%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "4"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04591"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.84%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "4"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006951"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-6.69%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "4"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1358"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.40%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "4"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002189"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.31%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "4"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009575"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.10%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "4"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006362"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.38%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "4"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.07%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "4"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.007975"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-19.32%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "4"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7469"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-8.98%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "4"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.07%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "4"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "4"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "4"
